# Update simulation results for the 380 kV case (pl_mw.xlsx, Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.6777082289699194
$ws.Range("C2").Value = 0.2002750294481057
$ws.Range("D2").Value = 0.04233067547223612
$ws.Range("F2").Value = 0.9286151925412298
$ws.Range("G2").Value = 0.00245660782892662
$ws.Range("K2").Value = 0.3217945786725807
$ws.Range("L2").Value = 0.3030788290682551
$ws.Range("N2").Value = 1.738431108733049
$ws.Range("O2").Value = 3.295122643933951

# Row 3
$ws.Range("B3").Value = 0.6333592719508658
$ws.Range("C3").Value = 0.2012658824147486
$ws.Range("D3").Value = 0.04044954733408446
$ws.Range("F3").Value = 0.9274437162160751
$ws.Range("G3").Value = 0.002459050399711322
$ws.Range("K3").Value = 0.2860315144321248
$ws.Range("L3").Value = 0.2919336572879843
$ws.Range("N3").Value = 1.756481030562854
$ws.Range("O3").Value = 3.304928512447987

# Row 4
$ws.Range("B4").Value = 0.6063794097092625
$ws.Range("C4").Value = 0.2019090694108527
$ws.Range("D4").Value = 0.03928133400210498
$ws.Range("F4").Value = 0.9272362135934458
$ws.Range("G4").Value = 0.002460631215894182
$ws.Range("K4").Value = 0.26409511240999
$ws.Range("L4").Value = 0.2852636740100536
$ws.Range("N4").Value = 1.768129488053266
$ws.Range("O4").Value = 3.312637939606248

# Row 5
$ws.Range("B5").Value = 0.5954485219349408
$ws.Range("C5").Value = 0.2021799446060903
$ws.Range("D5").Value = 0.03880198137560598
$ws.Range("F5").Value = 0.9272804059047672
$ws.Range("G5").Value = 0.002461295856489717
$ws.Range("K5").Value = 0.2551618970499163
$ws.Range("L5").Value = 0.2825892075816085
$ws.Range("N5").Value = 1.773018574599851
$ws.Range("O5").Value = 3.316204305763051

# Row 6
$ws.Range("B6").Value = 0.5936373172808942
$ws.Range("C6").Value = 0.2022254535997785
$ws.Range("D6").Value = 0.03872218678142048
$ws.Range("F6").Value = 0.9272955219808097
$ws.Range("G6").Value = 0.00246140745631269
$ws.Range("K6").Value = 0.2536789218443545
$ws.Range("L6").Value = 0.282147749757371
$ws.Range("N6").Value = 1.773838994650944
$ws.Range("O6").Value = 3.316822152806608

# Row 7
$ws.Range("B7").Value = 0.6062317335024545
$ws.Range("C7").Value = 0.2019126869855032
$ws.Range("D7").Value = 0.03927488260323742
$ws.Range("F7").Value = 0.9272362882012288
$ws.Range("G7").Value = 0.002460640096649377
$ws.Range("K7").Value = 0.2639746108677059
$ws.Range("L7").Value = 0.2852274285891809
$ws.Range("N7").Value = 1.768194848054655
$ws.Range("O7").Value = 3.312684317090714

# Row 8
$ws.Range("B8").Value = 0.6623650559994303
$ws.Range("C8").Value = 0.2006094654453534
$ws.Range("D8").Value = 0.04168481298356141
$ws.Range("F8").Value = 0.9281050776062969
$ws.Range("G8").Value = 0.002457433237843069
$ws.Range("K8").Value = 0.3094591110717033
$ws.Range("L8").Value = 0.2992000551065104
$ws.Range("N8").Value = 1.74453729818211
$ws.Range("O8").Value = 3.298153297963466

# Row 9
$ws.Range("B9").Value = 0.7744094459341113
$ws.Range("C9").Value = 0.1983290237220903
$ws.Range("D9").Value = 0.04630525738591018
$ws.Range("F9").Value = 0.9338684136866675
$ws.Range("G9").Value = 0.002451785091205074
$ws.Range("K9").Value = 0.3988154236440664
$ws.Range("L9").Value = 0.327973969643395
$ws.Range("N9").Value = 1.702633585672655
$ws.Range("O9").Value = 3.283054557363386

# Row 10
$ws.Range("B10").Value = 0.8579080837297397
$ws.Range("C10").Value = 0.1968200132249578
$ws.Range("D10").Value = 0.04963492100583267
$ws.Range("F10").Value = 0.9405780673657063
$ws.Range("G10").Value = 0.002448022016216954
$ws.Range("K10").Value = 0.4645497491867729
$ws.Range("L10").Value = 0.34995337302189
$ws.Range("N10").Value = 1.674582073930444
$ws.Range("O10").Value = 3.280130602906638

# Row 11
$ws.Range("B11").Value = 0.8961462873471078
$ws.Range("C11").Value = 0.1961693841487993
$ws.Range("D11").Value = 0.05113542769280599
$ws.Range("F11").Value = 0.9441681902478649
$ws.Range("G11").Value = 0.002446393229351705
$ws.Range("K11").Value = 0.4944696742279007
$ws.Range("L11").Value = 0.360135134519652
$ws.Range("N11").Value = 1.662414358936788
$ws.Range("O11").Value = 3.280575060483812

# Row 12
$ws.Range("B12").Value = 0.9106621541210984
$ws.Range("C12").Value = 0.1959281389110679
$ws.Range("D12").Value = 0.05170157402339726
$ws.Range("F12").Value = 0.9456050001101488
$ws.Range("G12").Value = 0.002445788331525148
$ws.Range("K12").Value = 0.5058016523861113
$ws.Range("L12").Value = 0.3640170330633907
$ws.Range("N12").Value = 1.657892064130439
$ws.Range("O12").Value = 3.280998541904722

# Row 13
$ws.Range("B13").Value = 0.9075343170362657
$ws.Range("C13").Value = 0.1959798673442457
$ws.Range("D13").Value = 0.05157973642844382
$ws.Range("F13").Value = 0.9452921196823638
$ws.Range("G13").Value = 0.002445918079415299
$ws.Range("K13").Value = 0.5033610292892945
$ws.Range("L13").Value = 0.3631798286328518
$ws.Range("N13").Value = 1.65886222385183
$ws.Range("O13").Value = 3.280895988748512

# Row 14
$ws.Range("B14").Value = 0.8973398010301139
$ws.Range("C14").Value = 0.1961494339691612
$ws.Range("D14").Value = 0.05118204637980028
$ws.Range("F14").Value = 0.9442848482993611
$ws.Range("G14").Value = 0.002446343226208503
$ws.Range("K14").Value = 0.4954019266376974
$ws.Range("L14").Value = 0.3604539741793644
$ws.Range("N14").Value = 1.662040595952629
$ws.Range("O14").Value = 3.280604786512129

# Row 15
$ws.Range("B15").Value = 0.8911000194399321
$ws.Range("C15").Value = 0.1962539664731473
$ws.Range("D15").Value = 0.05093818034445974
$ws.Range("F15").Value = 0.9436779322126085
$ws.Range("G15").Value = 0.002446605187487396
$ws.Range("K15").Value = 0.4905269857973451
$ws.Range("L15").Value = 0.3587877310781948
$ws.Range("N15").Value = 1.663998558583693
$ws.Range("O15").Value = 3.280459647824557

# Row 16
$ws.Range("B16").Value = 0.8554141859752633
$ws.Range("C16").Value = 0.1968632527565646
$ws.Range("D16").Value = 0.04953657202967321
$ws.Range("F16").Value = 0.9403542645643199
$ws.Range("G16").Value = 0.00244813012800019
$ws.Range("K16").Value = 0.462594712093761
$ws.Range("L16").Value = 0.3492916526844994
$ws.Range("N16").Value = 1.675389202812794
$ws.Range("O16").Value = 3.280137265009245

# Row 17
$ws.Range("B17").Value = 0.833586725435282
$ws.Range("C17").Value = 0.1972461934877359
$ws.Range("D17").Value = 0.04867308280586968
$ws.Range("F17").Value = 0.9384530372240505
$ws.Range("G17").Value = 0.002449086862925998
$ws.Range("K17").Value = 0.4454631783824539
$ws.Range("L17").Value = 0.3435130006590583
$ws.Range("N17").Value = 1.682528942193264
$ws.Range("O17").Value = 3.280394000786856

# Row 18
$ws.Range("B18").Value = 0.8210561193413639
$ws.Range("C18").Value = 0.1974698241772757
$ws.Range("D18").Value = 0.04817509400191256
$ws.Range("F18").Value = 0.9374101304272031
$ws.Range("G18").Value = 0.002449644972248664
$ws.Range("K18").Value = 0.4356112081053993
$ws.Range("L18").Value = 0.3402065221303445
$ws.Range("N18").Value = 1.68669134832393
$ws.Range("O18").Value = 3.280708700647551

# Row 19
$ws.Range("B19").Value = 0.8168176090983081
$ws.Range("C19").Value = 0.1975461216087098
$ws.Range("D19").Value = 0.04800625542280557
$ws.Range("F19").Value = 0.9370657167667034
$ws.Range("G19").Value = 0.002449835283463752
$ws.Range("K19").Value = 0.4322757962556807
$ws.Range("L19").Value = 0.3390899697188843
$ws.Range("N19").Value = 1.688110251923416
$ws.Range("O19").Value = 3.280843940784791

# Row 20
$ws.Range("B20").Value = 0.8359078197405552
$ws.Range("C20").Value = 0.1972050797999572
$ws.Range("D20").Value = 0.0487651408100831
$ws.Range("F20").Value = 0.9386501864817731
$ws.Range("G20").Value = 0.002448984207762092
$ws.Range("K20").Value = 0.4472866938164941
$ws.Range("L20").Value = 0.3441263632325899
$ws.Range("N20").Value = 1.681763127381709
$ws.Range("O20").Value = 3.280349383663321

# Row 21
$ws.Range("B21").Value = 0.9003332088495313
$ws.Range("C21").Value = 0.1960994889526049
$ws.Range("D21").Value = 0.05129891377342233
$ws.Range("F21").Value = 0.944578610492016
$ws.Range("G21").Value = 0.00244621802801504
$ws.Range("K21").Value = 0.4977396590308842
$ws.Range("L21").Value = 0.3612539107886761
$ws.Range("N21").Value = 1.661104713916576
$ws.Range("O21").Value = 3.280683394352423

# Row 22
$ws.Range("B22").Value = 0.9426478004000387
$ws.Range("C22").Value = 0.1954068373351454
$ws.Range("D22").Value = 0.05294284468298116
$ws.Range("F22").Value = 0.9489037724437992
$ws.Range("G22").Value = 0.002444479439986786
$ws.Range("K22").Value = 0.5307247331697909
$ws.Range("L22").Value = 0.3726009419406182
$ws.Range("N22").Value = 1.648100795602715
$ws.Range("O22").Value = 3.282389041040716

# Row 23
$ws.Range("B23").Value = 0.9200448168603543
$ws.Range("C23").Value = 0.1957737871435015
$ws.Range("D23").Value = 0.05206655732972365
$ws.Range("F23").Value = 0.9465541349619002
$ws.Range("G23").Value = 0.002445401036794948
$ws.Range("K23").Value = 0.5131191238851613
$ws.Range("L23").Value = 0.3665308213194436
$ws.Range("N23").Value = 1.654995680725987
$ws.Range("O23").Value = 3.281342615083105

# Row 24
$ws.Range("B24").Value = 0.8348583959786993
$ws.Range("C24").Value = 0.1972236564789895
$ws.Range("D24").Value = 0.04872352623512199
$ws.Range("F24").Value = 0.9385608991347496
$ws.Range("G24").Value = 0.002449030593025775
$ws.Range("K24").Value = 0.4464622912464904
$ws.Range("L24").Value = 0.3438490129656202
$ws.Range("N24").Value = 1.682109172606007
$ws.Range("O24").Value = 3.280369034557225

# Row 25
$ws.Range("B25").Value = 0.7438898985041931
$ws.Range("C25").Value = 0.1989166220521152
$ws.Range("D25").Value = 0.04506666404703452
$ws.Range("F25").Value = 0.9318746161228759
$ws.Range("G25").Value = 0.002453244894362175
$ws.Range("K25").Value = 0.3746263588055569
$ws.Range("L25").Value = 0.3200426166029899
$ws.Range("N25").Value = 1.773018838834911
$ws.Range("O25").Value = 3.285704748482601
